# a69_f10_bUPPachuca.xlsx - quarterly update (4to trimestre 2021 -> covers Oct-Dec 2021,
# reported in Jan 2022) + minor formatting cleanup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Selection moves to the "DESCRIPCIÓN" header block (G3:I3) ---
$ws.Range("G3:I3").Select() | Out-Null

# --- Row 3 header band: drop the oversized custom row height (back to default/auto) ---
$ws.Rows(3).AutoFit() | Out-Null

# --- Row 3: the merged label cells no longer wrap their (empty) text ---
$ws.Range("A3:C3").WrapText = $false

# --- Column M width tightens from 10.71 to a snug 8 (auto-fit to its content) ---
$ws.Columns(13).AutoFit() | Out-Null
$ws.Columns(13).ColumnWidth = 7.2

# --- Row 8 data: advance the reporting period from 3er to 4to trimestre 2021 ---
# Fecha de inicio / término del periodo que se informa
$ws.Range("B8").Value = 44470   # 2021-10-01
$ws.Range("C8").Value = 44561   # 2021-12-31

# Fecha de validación / Fecha de actualización
$ws.Range("K8").Value = 44571   # 2022-01-10
$ws.Range("L8").Value = 44571   # 2022-01-10

# Minor alignment tidy-up that rode along with the style cleanup in row 8
$ws.Range("D8").HorizontalAlignment = 1       # xlHAlignGeneral
$ws.Range("M8").HorizontalAlignment = -4131   # xlHAlignLeft
